# Sample Project / Main.xlsx - SAVE
# The rules table's "R20" row (row 9) has its "Integer max" value (column D)
# updated from 17 to 7.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D9").Value = 7
